$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestsRunner")

$ws.Columns("C").Insert()
$ws.Columns("C").NumberFormat = "@"

$data = @(
    @("TestCase","Description","version","Browser","Execute","username","password"),
    @("loginLogoutTest","validate OrangeHRM login and logout functionality","'latest","chrome","yes","admin","ad123"),
    @("loginLogoutTest","validate OrangeHRM login and logout functionality","'latest","chrome","yes","ad123","admin"),
    @("loginLogoutTest","validate OrangeHRM login and logout functionality","'latest","chrome","yes","Admin","admin123"),
    @("homePageTitleTest","validate title of home page","'latest","firefox","yes","Admin","admin123"),
    @("loginLogoutTest","validate OrangeHRM login and logout functionality","'95.0","chrome","yes","Admin","admin123"),
    @("homePageTitleTest","validate title of home page","'128.0","firefox","yes","Admin","admin123"),
    @("loginLogoutTest","validate OrangeHRM login and logout functionality","'latest","edge","yes","Admin","admin123"),
    @("homePageTitleTest","validate title of home page","'135.0","edge","yes","Admin","admin123")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}
